$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.586.80'
$ws.Cells.Item(2, 5).Value = '  +1.46%  '
$ws.Cells.Item(3, 4).Value = '1.791.78'
$ws.Cells.Item(3, 5).Value = '  -0.37%  '
$ws.Cells.Item(4, 4).Value = '''1.005'
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).Value = '''328.68'
$ws.Cells.Item(5, 5).Value = '  -2.56%  '
$ws.Cells.Item(6, 4).Value = '''0.9992'
$ws.Cells.Item(6, 5).Value = '  -0.20%  '
$ws.Cells.Item(7, 4).Value = '''0.4392'
$ws.Cells.Item(7, 5).Value = '  -3.24%  '
$ws.Cells.Item(8, 4).Value = '''0.3774'
$ws.Cells.Item(8, 5).Value = '  +6.38%  '
$ws.Cells.Item(9, 4).Value = '''45.69'
$ws.Cells.Item(9, 5).Value = '  +0.33%  '
$ws.Cells.Item(10, 4).Value = '''0.07627'
$ws.Cells.Item(10, 5).Value = '  +1.17%  '
$ws.Cells.Item(11, 4).Value = '''1.144'
$ws.Cells.Item(11, 5).Value = '  -0.23%  '
$ws.Cells.Item(12, 4).Value = '''22.72'
$ws.Cells.Item(12, 5).Value = '  +0.40%  '
$ws.Cells.Item(13, 4).Value = '''1.000'
$ws.Cells.Item(13, 5).Value = '  -0.35%  '
$ws.Cells.Item(14, 4).Value = '''6.280'
$ws.Cells.Item(14, 5).Value = '  +0.70%  '
$ws.Cells.Item(15, 4).Value = '''7.500'
$ws.Cells.Item(15, 5).Value = '  +3.29%  '
$ws.Cells.Item(16, 4).Value = '1.794.90'
$ws.Cells.Item(16, 5).Value = '  -0.19%  '
$ws.Cells.Item(17, 4).Value = '''0.00001094'
$ws.Cells.Item(17, 5).Value = '  +0.54%  '
$ws.Cells.Item(18, 4).Value = '''0.06711'
$ws.Cells.Item(18, 5).Value = '  +0.10%  '
$ws.Cells.Item(19, 4).Value = '''80.99'
$ws.Cells.Item(19, 5).Value = '  -0.57%  '
$ws.Cells.Item(20, 4).Value = '''1.001'
$ws.Cells.Item(20, 5).Value = '  -0.02%  '
$ws.Cells.Item(21, 4).Value = '''17.59'
$ws.Cells.Item(21, 5).Value = '  +2.07%  '
$ws.Cells.Item(22, 4).Value = '''6.253'
$ws.Cells.Item(22, 5).Value = '  -2.41%  '
$ws.Cells.Item(23, 4).Value = '28.615.19'
$ws.Cells.Item(23, 5).Value = '  +1.58%  '
$ws.Cells.Item(24, 4).Value = '''11.82'
$ws.Cells.Item(24, 5).Value = '  -1.06%  '
$ws.Cells.Item(25, 4).Value = '''2.445'
$ws.Cells.Item(25, 5).Value = '  +2.46%  '
$ws.Cells.Item(26, 4).Value = '''20.56'
$ws.Cells.Item(26, 5).Value = '  -0.35%  '
$ws.Cells.Item(27, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(27, 4).Value = '''2.396'
$ws.Cells.Item(27, 5).Value = '  -0.69%  '
$ws.Cells.Item(28, 2).Value = 'Monero'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(28, 4).Value = '''153.50'
$ws.Cells.Item(28, 5).Value = '  -1.02%  '
$ws.Cells.Item(29, 4).Value = '1.996.78'
$ws.Cells.Item(29, 5).Value = '  -0.34%  '
$ws.Cells.Item(30, 4).Value = '''1.339'
$ws.Cells.Item(30, 5).Value = '  +4.48%  '
$ws.Cells.Item(31, 4).Value = '''131.00'
$ws.Cells.Item(31, 5).Value = '  -1.96%  '
$ws.Cells.Item(32, 4).Value = '''3.972'
$ws.Cells.Item(32, 5).Value = '  -2.46%  '
$ws.Cells.Item(33, 4).Value = '''5.907'
$ws.Cells.Item(33, 5).Value = '  -0.12%  '
$ws.Cells.Item(34, 4).Value = '''0.09266'
$ws.Cells.Item(34, 5).Value = '  -2.62%  '
$ws.Cells.Item(35, 5).Value = '  +3.60%  '
$ws.Cells.Item(36, 4).Value = '''12.21'
$ws.Cells.Item(36, 5).Value = '  +0.52%  '
$ws.Cells.Item(37, 4).Value = '''0.06313'
$ws.Cells.Item(37, 5).Value = '  +0.88%  '
$ws.Cells.Item(38, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(38, 4).Value = '''5.242'
$ws.Cells.Item(38, 5).Value = '  +0.81%  '
$ws.Cells.Item(39, 2).Value = 'TheSandbox'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(39, 4).Value = '''0.6665'
$ws.Cells.Item(39, 5).Value = '  -0.76%  '
$ws.Cells.Item(40, 2).Value = 'VeChain'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(40, 4).Value = '''0.02336'
$ws.Cells.Item(40, 5).Value = '  -1.98%  '
$ws.Cells.Item(41, 4).Value = '''1.207'
$ws.Cells.Item(41, 5).Value = '  -0.64%  '
$ws.Cells.Item(42, 2).Value = 'WEMIXTOKEN'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(42, 4).Value = '''1.447'
$ws.Cells.Item(42, 5).Value = '  -2.55%  '
$ws.Cells.Item(43, 2).Value = 'FraxShare'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(43, 4).Value = '''8.073'
$ws.Cells.Item(43, 5).Value = '  -0.82%  '
$ws.Cells.Item(44, 4).Value = '''0.9989'
$ws.Cells.Item(44, 5).Value = '  -0.17%  '
$ws.Cells.Item(45, 4).Value = '''14.04'
$ws.Cells.Item(45, 5).Value = '  +0.25%  '
$ws.Cells.Item(46, 4).Value = '''0.6139'
$ws.Cells.Item(46, 5).Value = '  +0.54%  '
$ws.Cells.Item(47, 4).Value = '''3.818'
$ws.Cells.Item(47, 5).Value = '  -1.21%  '
$ws.Cells.Item(48, 4).Value = '''128.45'
$ws.Cells.Item(48, 5).Value = '  -0.73%  '
$ws.Cells.Item(49, 4).Value = '''2.031'
$ws.Cells.Item(49, 5).Value = '  -0.26%  '
$ws.Cells.Item(50, 4).Value = '''0.07015'
$ws.Cells.Item(50, 5).Value = '  -1.18%  '
$ws.Cells.Item(51, 4).Value = '''1.141'
$ws.Cells.Item(51, 5).Value = '  -2.17%  '
